$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (values first, then copy the formatting from H1,
# which already carries the bold/bordered/centered header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
